$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the current "IMF - Sales" column (column F)
# so that the new "IMF (20%)" columns land at F:G and the old IMF columns
# shift from F:G to H:I.
$ws.Range("F1:G1").EntireColumn.Insert() | Out-Null

# New header labels for the inserted columns
$ws.Range("F1").Value = "IMF (20%) - Sales"
$ws.Range("G1").Value = "IMF (20%) - Sales + Emp"

# New data values for the inserted columns, rows 2-8
$newF = @(-0.0365197584013312, 0.6334818562896358, 1.004697886772454, -0.304472963301639, 1.508446806846264, -8.108193661997705, 1.391211526531508)
$newG = @(0.03490700476505619, 0.5242771902968979, 2.052461775212154, -0.1838096993905943, 1.052810667356595, 15.86246084502586, 1.603193569135634)

for ($i = 0; $i -lt 7; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $newF[$i]
    $ws.Cells.Item($row, 7).Value = $newG[$i]
}

# The previously existing "OECD (20%) - Sales" / "OECD (20%) - Sales + Emp"
# columns (now at J:K after the insert, since the old F:G IMF columns shifted
# to H:I) are removed - that data no longer belongs in the sheet.
$ws.Range("J1:K1").EntireColumn.Delete() | Out-Null

$wb.Save()
